# Update the city-tax settlement status sheet from 令和4年度 (FY R4) to
# 令和5年度 (FY R5): refresh the year label and all the settlement figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B: fiscal year label for every data row (2-22)
$ws.Range("B2:B22").Value = "令和5年度"

# Row-by-row refreshed figures (F adjustedAmount, G incomeAmount,
# H nonPaymentDeficit, I unpaidIncome, J incomeRate, K comparedToThePreviousYear)
$data = @(
    @(2,  24699660604, 24502700421, 1917876,   202470148, 99.181089969930483, 101.4),
    @(3,  549194597,   186574140,   28903205,  333945870, 34,                  99.4),
    @(4,  6572491300,  6545925632,  160000,    26405668,  99.6,                94.1),
    @(5,  58868576,    10738014,    4729309,   43551553,  18.2,                76.1),
    @(6,  23799769900, 23601407402, 5598500,   193988898, 99.2,               101.9),
    @(7,  729318186,   187750389,   78750568,  462911218, 25.7,                96.5),
    @(8,  3435696600,  3433850000, 0,          2200100,   99.939605996079123, 98.5),
    @(9,  8349000,     2723742,     721400,    4903858,   32.6,                63.3),
    @(10, 258988200,   258988200,  0,          0,         100,                100.4),
    @(11, 1328767100,  1310135633,  70800,     18664639,  98.6,               103.1),
    @(12, 57935400,    57935400,   0,          0,         100,                 98.7),
    @(13, 59000321,    14527831,    5851257,   38628105,  24.6,                90.9),
    @(14, 2997969112,  2997969112, 0,          0,         100,                100.6),
    @(16, 30957300,    30957300,   0,          0,         100,                116.5),
    @(18, 2346296600,  2343742400, 0,          2554200,   99.9,                98.5),
    @(19, 33460900,    5254100,     9404200,   18802600,  15.7,                33.6),
    @(20, 65528532116, 65083611500, 7747176,   446283653, 99.275801897828273, 100.6),
    @(21, 1438199230,  407568216,   128367589, 902743204, 28.3,                94.3),
    @(22, 66966731346, 65491179716, 136114765, 1349026857,97.8,               100.5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
    $ws.Cells.Item($r, 8).Value = $row[3]
    $ws.Cells.Item($r, 9).Value = $row[4]
    $ws.Cells.Item($r, 10).Value = $row[5]
    $ws.Cells.Item($r, 11).Value = $row[6]
}

# Row 15: all zero values, rate columns show "-" (text)
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = "-"
$ws.Cells.Item(15, 11).Value = "-"

# Row 17: unpaid income moved to H, I cleared to 0, J stays 0, K becomes "-"
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 7650
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = "-"

# Cosmetic: last selection, matching the refreshed workbook
$ws.Range("L21").Select()
